$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New ShipmentTrackNum values for rows 2..22 (column C), mirrored into column D
# for the rows where D previously held the same value as C.
$newValues = @(
    "320018586090",
    "320018586104",
    "320018586137",
    "320018586159",
    "320018586192",
    "320018586218",
    "320018586240",
    "320018586262",
    "320018586295",
    "320018586310",
    "320018586354",
    "320018586376",
    "320018586402",
    "320018586424",
    "320018586457",
    "320018586479",
    "320018586516",
    "320018586538",
    "320018586560",
    "320018586582",
    "320018586619"
)

# Rows where column D mirrors column C's value.
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $val = $newValues[$i]
    $ws.Cells.Item($row, 3).Value = $val
    if ($mirrorRows -contains $row) {
        $ws.Cells.Item($row, 4).Value = $val
    }
}
